$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = 55
$ws1.Range("F6").Value = 7406
$ws1.Range("F9").Value = 6639
$ws1.Range("F11").Value = 283
$ws1.Range("F12").Value = 4636
$ws1.Range("F16").Value = 4791
$ws1.Range("F17").Value = 4791
$ws1.Range("F18").Value = 1043
$ws1.Range("F19").Value = 269
$ws1.Range("F22").Value = 375
$ws1.Range("F23").Value = 291
$ws1.Range("F26").Value = 179
$ws1.Range("F30").Value = 80
$ws1.Range("F31").Value = 8452
$ws1.Range("F33").Value = 1478
$ws1.Range("F37").Value = 130
$ws1.Range("F40").Value = 1741
$ws1.Range("F41").Value = 1011
$ws1.Range("F42").Value = 47
$ws1.Range("F43").Value = 4407
$ws1.Range("F44").Value = 360
$ws1.Range("F45").Value = 124
$ws1.Range("F47").Value = 861
$ws1.Range("F48").Value = 1155

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 24

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G4").Value = 55
$ws4.Range("F8").Value = 7406
$ws4.Range("F11").Value = 6639
$ws4.Range("F13").Value = 283
$ws4.Range("F15").Value = 4636
$ws4.Range("F19").Value = 4791
$ws4.Range("F20").Value = 1043
$ws4.Range("F21").Value = 269
$ws4.Range("F23").Value = 375
$ws4.Range("F24").Value = 291
$ws4.Range("F27").Value = 179
$ws4.Range("F31").Value = 80
$ws4.Range("F33").Value = 8452
$ws4.Range("F35").Value = 1478
$ws4.Range("F39").Value = 130
$ws4.Range("F42").Value = 1741
$ws4.Range("F43").Value = 1011
$ws4.Range("F44").Value = 47
$ws4.Range("F45").Value = 360
$ws4.Range("F46").Value = 124
$ws4.Range("F48").Value = 861
$ws4.Range("F49").Value = 1155
